# Generate Report for Handoff
#
# The 84d8f86b-6d96-4774-8943-d34f44e81218 file has been handed off / is no
# longer pending, so its row is removed from every sheet (Overview, zh-cn,
# de-de). The remaining 293fed2a-... entry moves from "Handed back: in sync
# with en-US" to "Ready for handoff", and its handoff timestamps are bumped
# to reflect the new handoff.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Rows.Item(3).Delete()
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

# --- zh-cn sheet ------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Rows.Item(3).Delete()
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "2016-01-29 02:29:40"

# --- de-de sheet ------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Rows.Item(3).Delete()
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "2016-01-29 02:29:51"
